$wb = $excel.ActiveWorkbook

# Use sheet "pubmed_fulltext1" as a template and duplicate it, placing the
# copy immediately after the template (as the new last sheet).
$template = $wb.Worksheets.Item("pubmed_fulltext1")
$template.Copy($null, $template)

# The copied sheet becomes the active sheet
$ws = $wb.ActiveSheet
$ws.Name = "pubmed_fulltext2"

# Header row (values are identical to the template, but set explicitly for clarity)
$ws.Range("B1").Value = "Fold"
$ws.Range("C1").Value = "Version"
$ws.Range("D1").Value = "Epoch"
$ws.Range("E1").Value = "Recall"
$ws.Range("F1").Value = "Precision"
$ws.Range("G1").Value = "Accuracy"
$ws.Range("H1").Value = "Fbeta"
$ws.Range("I1").Value = "Best Recall"
$ws.Range("J1").Value = "Best Precision"
$ws.Range("K1").Value = "Best Threshold"
$ws.Range("L1").Value = "False Neg(0.5)"
$ws.Range("M1").Value = "False Pos(0.5)"
$ws.Range("N1").Value = "Val loss"

# Columns E:H hold numeric-looking values that are stored as text in the
# source workbook. Force Text number format before assigning them so the
# values are kept as strings, then restore the default style so the cells
# do not end up with an extra explicit style index.
$textRange = $ws.Range("E2:H5")
$textRange.NumberFormat = "@"

# Row 2 - fold_0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "fold_0"
$ws.Range("C2").Value = "13.10_09.26"
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = "0.8467742"
$ws.Range("F2").Value = "0.41501975"
$ws.Range("G2").Value = "0.86510503"
$ws.Range("H2").Value = "0.7009346"
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0.5
$ws.Range("L2").Value = 19
$ws.Range("M2").Value = 148
$ws.Range("N2").Value = 0.6406238003964385

# Row 3 - fold_1
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "fold_1"
$ws.Range("C3").Value = "13.10_09.41"
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = "0.92"
$ws.Range("F3").Value = "0.42910448"
$ws.Range("G3").Value = "0.868336"
$ws.Range("H3").Value = "0.74869794"
$ws.Range("I3").Value = 0.968
$ws.Range("J3").Value = 0.4158
$ws.Range("K3").Value = 0.2996
$ws.Range("L3").Value = 10
$ws.Range("M3").Value = 153
$ws.Range("N3").Value = 0.7133551667774877

# Row 4 - fold_2
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "fold_2"
$ws.Range("C4").Value = "13.10_09.58"
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = "0.928"
$ws.Range("F4").Value = "0.4566929"
$ws.Range("G4").Value = "0.8812601"
$ws.Range("H4").Value = "0.7692308"
$ws.Range("I4").Value = 0.976
$ws.Range("J4").Value = 0.404
$ws.Range("K4").Value = 0.1213
$ws.Range("L4").Value = 9
$ws.Range("M4").Value = 138
$ws.Range("N4").Value = 0.5740149182177359

# Row 5 - fold_3
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "fold_3"
$ws.Range("C5").Value = "13.10_10.16"
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = "0.944"
$ws.Range("F5").Value = "0.43223444"
$ws.Range("G5").Value = "0.8691438"
$ws.Range("H5").Value = "0.76326"
$ws.Range("I5").Value = 0.96
$ws.Range("J5").Value = 0.4364
$ws.Range("K5").Value = 0.4906
$ws.Range("L5").Value = 7
$ws.Range("M5").Value = 155
$ws.Range("N5").Value = 0.5730994257234758

# Restore the default (unstyled) format on the text range now that the
# string values have been written.
$textRange.Style = "Normal"

# Restore the originally active sheet/selection so the workbook-level view
# state is left unchanged by this edit.
$wb.Worksheets.Item("medbert").Activate()
